$wb = $excel.ActiveWorkbook

# =========================================================================
# Sheet "ITR input data" -- add new row 33 (PPL Corp.)
# =========================================================================
$wsITR = $wb.Worksheets.Item("ITR input data")

$fmtCols4 = @("B","G","I","O","P","AE","AF","AG","AH","AI")
foreach ($col in $fmtCols4) {
    $src = $col + "32"
    $dst = $col + "33"
    $wsITR.Range($src).Copy()
    $wsITR.Range($dst).PasteSpecial(-4122)
}

$wsITR.Range("A33").Value = "PPL Corp."
$wsITR.Range("B33").Value = "9N3UAJSNOUXFKQLF3V18"
$wsITR.Range("C33").Value = "US69351T1060"
$wsITR.Range("D33").Value = "US"
$wsITR.Range("E33").Value = "North America"
$wsITR.Range("F33").Value = "Electricity Utilities"
$wsITR.Range("G33").Value = "equity"
$wsITR.Range("H33").Value = "USD"
$wsITR.Range("I33").Value = 44196
$wsITR.Range("J33").Value = 19865342074
$wsITR.Range("K33").Value = 7769000000
$wsITR.Range("L33").Value = 40943342074
$wsITR.Range("M33").Value = 41758342074
$wsITR.Range("N33").Value = 45680000000
$wsITR.Range("O33").Value = "Mt CO2"
$wsITR.Range("P33").Value = "TWh"
$wsITR.Range("AE33").Value = 30.08848722
$wsITR.Range("AF33").Value = 30.24837145
$wsITR.Range("AG33").Value = 31.61146904
$wsITR.Range("AH33").Value = 28.77891532
$wsITR.Range("AI33").Value = 28.07780713
$wsITR.Range("AS33").Value = 38.35525864
$wsITR.Range("AT33").Value = 37.44283235
$wsITR.Range("AU33").Value = 39.59007518
$wsITR.Range("AV33").Value = 35.15293172
$wsITR.Range("AW33").Value = 32.487984334642732

# View: unfreeze/refreeze so the pane's top-left moves off row/col 1, then
# select the new row's first few cells.
$wsITR.Activate()
$excel.ActiveWindow.FreezePanes = $false
$wsITR.Range("B2").Select()
$excel.ActiveWindow.FreezePanes = $true
$wsITR.Range("A33:C33").Select()

# =========================================================================
# Sheet "ITR target input data" -- add new rows 42 & 43 (PPL Corp. targets)
# =========================================================================
$wsTgt = $wb.Worksheets.Item("ITR target input data")

$fmtCols5 = @("A","B","C","D","E","F","G","H","J","K","L")
foreach ($col in $fmtCols5) {
    $src = $col + "40"
    $dst42 = $col + "42"
    $dst43 = $col + "43"
    $wsTgt.Range($src).Copy()
    $wsTgt.Range($dst42).PasteSpecial(-4122)
    $wsTgt.Range($src).Copy()
    $wsTgt.Range($dst43).PasteSpecial(-4122)
}

$wsTgt.Range("A42").Value = "PPL Corp."
$wsTgt.Range("B42").Value = "9N3UAJSNOUXFKQLF3V18"
$wsTgt.Range("C42").Value = "US69351T1060"
$wsTgt.Range("D42").Value = 2050
$wsTgt.Range("E42").Value = "absolute"
$wsTgt.Range("F42").Value = "S1+S2"
$wsTgt.Range("G42").Value = 2021
$wsTgt.Range("H42").Value = 2010
$wsTgt.Range("I42").Formula = "=60736086+1597157"
$wsTgt.Range("J42").Value = "t CO2"
$wsTgt.Range("K42").Value = 2035
$wsTgt.Range("L42").Value = 0.7

$wsTgt.Range("A43").Value = "PPL Corp."
$wsTgt.Range("B43").Value = "9N3UAJSNOUXFKQLF3V18"
$wsTgt.Range("C43").Value = "US69351T1060"
$wsTgt.Range("D43").Value = 2050
$wsTgt.Range("E43").Value = "absolute"
$wsTgt.Range("F43").Value = "S1+S2"
$wsTgt.Range("G43").Value = 2021
$wsTgt.Range("H43").Value = 2010
$wsTgt.Range("I43").Formula = "=60736086+1597157"
$wsTgt.Range("J43").Value = "t CO2"
$wsTgt.Range("K43").Value = 2040
$wsTgt.Range("L43").Value = 0.8

$wsTgt.Activate()
$excel.ActiveWindow.FreezePanes = $false
$wsTgt.Range("B9").Select()
$excel.ActiveWindow.FreezePanes = $true
$wsTgt.Range("A43").Select()

# =========================================================================
# Sheet "ITR input data (2)" -- view only (no data change)
# =========================================================================
$wsITR2 = $wb.Worksheets.Item("ITR input data (2)")
$wsITR2.Activate()
$excel.ActiveWindow.FreezePanes = $false
$wsITR2.Range("L14").Select()
$excel.ActiveWindow.FreezePanes = $true
$wsITR2.Range("A40:XFD40").Select()

# =========================================================================
# Sheet "ITR target input data (2)" -- view only (no data change)
# =========================================================================
$wsTgt2 = $wb.Worksheets.Item("ITR target input data (2)")
$wsTgt2.Activate()
$excel.ActiveWindow.TopLeftCell = $wsTgt2.Range("A13")

# =========================================================================
# Sheet "Portfolio" -- add new row 33 (PPL Corp.)
# =========================================================================
$wsPort = $wb.Worksheets.Item("Portfolio")

$fmtCols7 = @("A","B","C","D")
foreach ($col in $fmtCols7) {
    $src = $col + "32"
    $dst = $col + "33"
    $wsPort.Range($src).Copy()
    $wsPort.Range($dst).PasteSpecial(-4122)
}

$wsPort.Range("A33").Value = "PPL Corp."
$wsPort.Range("B33").Value = "9N3UAJSNOUXFKQLF3V18"
$wsPort.Range("C33").Value = "US69351T1060"
$wsPort.Range("D33").Value = "US69351T1060"
$wsPort.Range("E33").Formula = "=RANDBETWEEN(35000,250000)"

# Portfolio becomes the active sheet/tab, with E33 selected -- activate it
# LAST so it ends up as the workbook's active tab.
$wsPort.Activate()
$wsPort.Range("E33").Select()
